$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "0.9993" or "241.02"
# are not auto-converted to numbers by Excel, matching the original inlineStr formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.251.00"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.863.55"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "0.7142"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "241.02"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D8").Value = "0.07737"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").Value = "0.3085"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").Value = "24.94"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "0.08342"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").Value = "1.883.63"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("D14").Value = "0.7147"
$ws.Range("E14").Value = "  -2.05%  "

$ws.Range("D15").Value = "90.99"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").Value = "29.269.71"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "5.955"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "242.95"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").Value = "0.000007824"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "2.131.54"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "7.900"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "0.1605"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").Value = "163.36"

$ws.Range("D27").Value = "8.903"
$ws.Range("E27").Value = "  -1.70%  "

$ws.Range("D28").Value = "18.56"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "1.345"
$ws.Range("E29").Value = "  -1.66%  "

$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").Value = "4.422"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").Value = "4.261"
$ws.Range("E32").Value = "  +2.54%  "

$ws.Range("D33").Value = "0.8615"
$ws.Range("E33").Value = "  +19.12%  "

$ws.Range("D34").Value = "0.05152"
$ws.Range("E34").Value = "  -2.50%  "

$ws.Range("D35").Value = "1.930"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -2.75%  "

$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("D39").Value = "2.690"
$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("D40").Value = "1.168.83"
$ws.Range("E40").Value = "  -5.36%  "

$ws.Range("D41").Value = "6.203"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").Value = "0.8973"
$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").Value = "72.92"
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "102.13"
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("D46").Value = "2.025.53"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.5183"
$ws.Range("E47").Value = "  -2.85%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.790"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.326"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").Value = "7.061"
$ws.Range("E51").Value = "  -0.49%  "

# Restore default cell style on column D so no stray number-format style is left behind.
$ws.Range("D2:D51").Style = "Normal"
